$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 17..64 hold "Periodo Mora" (col E) + "Valor Mora" (col F) for the same
# worker, originally listed newest-to-oldest (2009 down to 1610, with a
# duplicate/garbled 1610 left over at row 64). The update sorts this block
# into chronological (ascending) order, so row 17 becomes the new duplicate
# 1610 entry and row 64 becomes 2009.
$periods = @(17,18,19,20,21,22,23,24,25,26,27,28,29,30,31,32,33,34,35,36,37,38,39,40,41,42,43,44,45,46,47,48,49,50,51,52,53,54,55,56,57,58,59,60,61,62,63,64)

$eVals = @()
$fVals = @()
foreach ($r in $periods) {
    $eVals += , ($ws.Range("E$r").Value2)
    $fVals += , ($ws.Range("F$r").Value2)
}

$n = $periods.Length
for ($i = 0; $i -lt $n; $i++) {
    $r = $periods[$i]
    $ws.Range("E$r").Value2 = $eVals[$n - 1 - $i]
    $ws.Range("F$r").Value2 = $fVals[$n - 1 - $i]
}

# Column C ("N° Doc Trabajador") width was manually narrowed.
$ws.Columns("C").ColumnWidth = 8

# The logo picture shifts left (column C got narrower), keeping its size.
$logo = $ws.Shapes.Item(1)
$logo.Left = $logo.Left - 19.0
